$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.446.96'
$ws.Range('E2').Value = '  +3.94%  '
$ws.Range('D3').Value = '2.430.45'
$ws.Range('E3').Value = '  +2.94%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.571'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.95%  '
$ws.Range('E9').Value = '  +4.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.81'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.38%  '
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('E12').Value = '  -1.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.90'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.99%  '
$ws.Range('D14').Value = '2.863.96'
$ws.Range('E14').Value = '  +2.95%  '
$ws.Range('D15').Value = '60.349.03'
$ws.Range('E15').Value = '  +3.85%  '
$ws.Range('E16').Value = '  +3.66%  '
$ws.Range('D17').Value = '2.430.01'
$ws.Range('E17').Value = '  +3.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.43'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.41'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '335.49'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.87'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.99'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.172'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.63'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('D28').Value = '0.0₃0785'
$ws.Range('E28').Value = '  +5.94%  '
$ws.Range('E29').Value = '  +2.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.29'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '169.50'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.82'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.38%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.31'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.25'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.63'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '40.11'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.07%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.418'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.19%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '317.94'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.43%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.71'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.53%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '142.75'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.02%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.90'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.32%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0525'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.89%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0958'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.28%  '
$ws.Range('B47').Value = 'Polygon'
$ws.Range('C47').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.408'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.84%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.573'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('E49').Value = '  +0.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.78%  '
$ws.Range('E51').Value = '  -0.12%  '
